# Apply the crypto-price / 1h-volume-change refresh described by the commit
# "Updated cryptos list ... with GitHub Actions": a batch of D (Price) and
# E (Volume 1h) text-cell updates on the single data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price column (D) ---------------------------------------------------
# The Price column stores plain text (e.g. "245.83", "37.127.37") rather than
# numbers. Assigning a purely-numeric-looking string like "245.88" through
# .Value would normally be auto-coerced into a numeric cell, so we briefly
# force the whole Price column to Text format, write the new values, then
# restore the default (unstyled) look -- this keeps every Price cell a text
# cell exactly like the rest of the (untouched) rows, with no visible format
# change.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "37.202.88"
$ws.Range("D3").Value = "2.001.54"
$ws.Range("D5").Value = "245.88"
$ws.Range("D7").Value = "59.94"
$ws.Range("D9").Value = "0.385"
$ws.Range("D10").Value = "0.0806"
$ws.Range("D12").Value = "15.05"
$ws.Range("D13").Value = "22.47"
$ws.Range("D14").Value = "2.293.14"
$ws.Range("D16").Value = "5.43"
$ws.Range("D17").Value = "2.001.04"
$ws.Range("D18").Value = "37.125.72"
$ws.Range("D19").Value = "70.17"
$ws.Range("D20").Value = "0.0₃0864"
$ws.Range("D21").Value = "5.18"
$ws.Range("D22").Value = "230.04"
$ws.Range("D24").Value = "2.47"
$ws.Range("D26").Value = "9.42"
$ws.Range("D27").Value = "0.142"
$ws.Range("D28").Value = "163.88"
$ws.Range("D29").Value = "19.64"
$ws.Range("D30").Value = "1.36"
$ws.Range("D32").Value = "4.80"
$ws.Range("D33").Value = "0.0656"
$ws.Range("D34").Value = "4.49"
$ws.Range("D37").Value = "1.80"
$ws.Range("D38").Value = "3.30"
$ws.Range("D39").Value = "5.38"
$ws.Range("D40").Value = "0.0982"
$ws.Range("D44").Value = "16.68"
$ws.Range("D45").Value = "91.17"
$ws.Range("D46").Value = "1.369.63"
$ws.Range("D48").Value = "7.44"
$ws.Range("D49").Value = "2.05"
$ws.Range("D50").Value = "2.84"
$ws.Range("D51").Value = "46.33"

$priceRange.Style = "Normal"

# --- Volume(1h) column (E) -----------------------------------------------
# These cells already contain padding spaces and a trailing "%", so they stay
# text automatically; no format juggling required.
$ws.Range("E2").Value = "  +1.50%  "
$ws.Range("E3").Value = "  +2.15%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("E5").Value = "  +0.60%  "
$ws.Range("E6").Value = "  +1.56%  "
$ws.Range("E7").Value = "  +1.23%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  +2.74%  "
$ws.Range("E10").Value = "  +2.33%  "
$ws.Range("E11").Value = "  +1.10%  "
$ws.Range("E12").Value = "  +5.22%  "
$ws.Range("E13").Value = "  +4.74%  "
$ws.Range("E14").Value = "  +2.13%  "
$ws.Range("E15").Value = "  +0.21%  "
$ws.Range("E16").Value = "  +2.54%  "
$ws.Range("E17").Value = "  +2.37%  "
$ws.Range("E18").Value = "  +1.51%  "
$ws.Range("E19").Value = "  +0.61%  "
$ws.Range("E20").Value = "  +1.49%  "
$ws.Range("E21").Value = "  +2.23%  "
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("E24").Value = "  +0.82%  "
$ws.Range("E25").Value = "  +0.12%  "
$ws.Range("E26").Value = "  +2.13%  "
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("E28").Value = "  +1.90%  "
$ws.Range("E29").Value = "  +0.92%  "
$ws.Range("E30").Value = "  +11.46%  "
$ws.Range("E31").Value = "  +0.93%  "
$ws.Range("E32").Value = "  +0.93%  "
$ws.Range("E33").Value = "  +6.54%  "
$ws.Range("E34").Value = "  +2.04%  "
$ws.Range("E35").Value = "  +3.65%  "
$ws.Range("E36").Value = "  +0.19%  "
$ws.Range("E37").Value = "  +1.90%  "
$ws.Range("E38").Value = "  -6.47%  "
$ws.Range("E39").Value = "  -2.11%  "
$ws.Range("E40").Value = "  +0.07%  "
$ws.Range("E41").Value = "  +0.67%  "
$ws.Range("E42").Value = "  +1.77%  "
$ws.Range("E43").Value = "  +0.63%  "
$ws.Range("E44").Value = "  +5.25%  "
$ws.Range("E45").Value = "  +3.26%  "
$ws.Range("E46").Value = "  -0.06%  "
$ws.Range("E48").Value = "  +4.07%  "
$ws.Range("E49").Value = "  +11.70%  "
$ws.Range("E50").Value = "  +0.00%  "
$ws.Range("E51").Value = "  +4.72%  "
